# Add rule for draw down and top up of undrawn balance.
# The draw down is based on the CCF column (and can be overridden by scenario
# input); the top up can be specified in scenario input as well.
#
# Sheet "intangible redemption" currently has 4 columns (A:D):
#   A = Date, B = Repaymentrate (CCF%), C = Repayment, D = Prepayment
# We insert two new scenario-input columns right after B (percentages, like
# B) and two new amount columns at the end (same style as the old C/D
# amount columns), shifting the existing Repayment/Prepayment amount
# columns from C/D to E/F:
#   A = Date
#   B = Repaymentrate
#   C = TopUpRate        (new, percentage input)
#   D = DrawDownRate      (new, percentage input)
#   E = Repayment         (was C)
#   F = Prepayment        (was D)
#   G = TopUp              (new, amount)
#   H = DrawDown           (new, amount)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intangible redemption")

# --- Move the existing Repayment/Prepayment amount columns from C:D to E:F ---
# (done as values-then-formats; a single xlPasteAll paste does not carry the
# cell style index through in this host). Row 1 has no data in C/D so it is
# left untouched.
$ws.Range("C2:D25").Copy()
$ws.Range("E2:F25").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("C2:D25").Copy()
$ws.Range("E2:F25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 category headers ---
$ws.Range("C2").Value = "Intangible assets"
$ws.Range("D2").Value = "Intangible assets"
$ws.Range("G2").Value = "Other loans"
$ws.Range("H2").Value = "Other loans"

# --- Row 3 metric headers (write in this order so new shared strings land
#     at the same indices Excel would have produced) ---
$ws.Range("G3").Value = "TopUp"
$ws.Range("H3").Value = "DrawDown"
$ws.Range("C3").Value = "TopUpRate"
$ws.Range("D3").Value = "DrawDownRate"

# --- New rate columns C:D get the same percentage formatting as column B ---
$ws.Range("B4:B25").Copy()
$ws.Range("C4:D25").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C4:D25").Value = 0

# --- New amount columns G:H start out unstyled, same as the majority of
#     cells in the amount columns they are modelled on (E4, a cell with no
#     override, carries the unstyled baseline - the few cells that do carry
#     the "Comma" style are overrides, handled individually below) ---
$ws.Range("G4:H25").Value = 0

$ws.Range("E4").Copy()
$ws.Range("G4:H25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Scenario overrides ---
# Row 6 (2025-10-31): 10% draw-down rate override, drawing down 100,000
$ws.Range("C6").Value = 0.1
$ws.Range("G6").Value = 100000

# Row 10 (2025-08-29): 50% top-up rate override
$ws.Range("D10").Value = 0.5

# Row 11 (2025-09-29): top-up of 100,000
$ws.Range("H11").Value = 100000

# Non-zero amount overrides use the same "Comma" number style as the other
# populated amount cells (e.g. E7)
$ws.Range("E7").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("E7").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths: widen the two new rate columns (same width as column
#     B) and keep the amount columns at the original "bestFit" width ---
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 22.140625
$ws.Range("E1:H1").EntireColumn.ColumnWidth = 11.140625

# --- Selection, matching the post-edit workbook ---
$ws.Range("D11").Select()
